$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1917
$ws1.Range("F6").Value = 2810
$ws1.Range("F21").Value = 17
$ws1.Range("F22").Value = 3
$ws1.Range("F23").Value = 21
$ws1.Range("F25").Value = 26
$ws1.Range("F27").Value = 1794
$ws1.Range("F29").Value = 430
$ws1.Range("F33").Value = 317

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1917
$ws4.Range("F7").Value = 2810
$ws4.Range("F22").Value = 17
$ws4.Range("F23").Value = 3
$ws4.Range("F24").Value = 21
$ws4.Range("F26").Value = 26
$ws4.Range("F28").Value = 1794
$ws4.Range("F30").Value = 430
$ws4.Range("F34").Value = 317
